$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (shifts existing rows 8-15 down to 9-16),
# introducing the new "discount_upto" field row right after "discount_value".
$ws.Rows("8:8").Insert()

# Populate the newly inserted row with the discount_upto field definition.
$ws.Range("A8").Value = "discount_upto"
$ws.Range("B8").Value = "Yes"
$ws.Range("C8").Value = "number"
$ws.Range("K8").Value = 100.0
